$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 / J1 - copy style/format from existing header cell (H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I and J, rows 2-38
$data = @(
    @(1,3),
    @(1,4),
    @(1,6),
    @(1,6),
    @(1,4),
    @(1,6),
    @(1,6),
    @(1,6),
    @(1,5),
    @(1,6),
    @(1,6),
    @(1,6),
    @(1,3),
    @(1,5),
    @(1,5),
    @(1,6),
    @(1,5),
    @(1,7),
    @(1,6),
    @(1,6),
    @(1,5),
    @(1,6),
    @(1,6),
    @(1,7),
    @(1,6),
    @(1,6),
    @(1,5),
    @(1,3),
    @(1,6),
    @(1,6),
    @(1,5),
    @(7,9),
    @(1,3),
    @(6,9),
    @(1,5),
    @(3,5),
    @(1,2)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
